$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (56) with the latest run's date + allocation split,
# matching the style of the existing data rows (plain values, default style).
$ws.Range("A56").Value = "'10/27/2025"
$ws.Range("A56").Style = "Normal"
$ws.Range("B56").Value = 0.189372928692943
$ws.Range("C56").Value = 0.810627071307057
